# "Merge ScanMe card back into mobster pitch"
# The "DWK" card/account is renamed to "Rocko" throughout, the stale
# "Isaac" contact is replaced with "Dan Velvet", and a new "Ponzi"
# contact/sale (paid via the same $DorkyDomains / CashApp rail as Rocko)
# is appended to both the "User" and "Sales" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "User": contact/payment directory
# ---------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("User")

# Row 2: "DWK" -> "Rocko" (same Cash App / $DorkyDomains / cell number)
$wsUser.Range("A2").Value = "Rocko"

# Column A now holds ID text, so size it to fit like the rest of the sheet
$wsUser.Columns.Item(1).AutoFit()

# Row 5: "Isaac" -> "Dan Velvet" (cell number untouched)
$wsUser.Range("A5").Value = "Dan Velvet"

# New row 6: "Ponzi" joins on the same $DorkyDomains Cash App rail
$wsUser.Range("A6").Value = "Ponzi"
$wsUser.Range("B6").Value = "CashApp"
$wsUser.Range("C6").Value = '$DorkyDomains'
$wsUser.Range("D6").Value = "503-975-5544"

# Leave the selection where the new row was typed
[void]$wsUser.Range("C6:D6").Select()

# ---------------------------------------------------------------
# Sheet "Sales": sale log
# ---------------------------------------------------------------
$wsSales = $wb.Worksheets.Item("Sales")

# Existing DWK sales now belong to Rocko
$wsSales.Range("A2").Value = "Rocko"
$wsSales.Range("D2").Value = "Rocko"
$wsSales.Range("A3").Value = "Rocko"
$wsSales.Range("D3").Value = "Rocko"
$wsSales.Range("A4").Value = "Rocko"
$wsSales.Range("D4").Value = "Rocko"

# New row 5: Ponzi's sale, same free/unpaid status
$wsSales.Range("A5").Value = "Rocko"
$wsSales.Range("B5").Value = 4
$wsSales.Range("C5").Value = "Ponzi"
$wsSales.Range("D5").Value = "Rocko"
$wsSales.Range("E5").Value = "Free"

[void]$wsSales.Range("F5").Select()
